$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text so numeric-looking values are not
# auto-converted to floating point numbers by Excel (matches original inlineStr text).
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '42.756.48'
$ws.Range('E2').Value = '  +1.05%  '

$ws.Range('D3').Value = '2.522.71'
$ws.Range('E3').Value = '  +0.17%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '315.98'
$ws.Range('E5').Value = '  +4.28%  '

$ws.Range('D6').Value = '94.57'
$ws.Range('E6').Value = '  -1.46%  '

$ws.Range('D7').Value = '0.578'
$ws.Range('E7').Value = '  -0.81%  '

$ws.Range('E8').Value = '  -0.16%  '

$ws.Range('D9').Value = '0.528'
$ws.Range('E9').Value = '  -1.27%  '

$ws.Range('D10').Value = '35.83'
$ws.Range('E10').Value = '  -1.48%  '

$ws.Range('D11').Value = '0.0809'
$ws.Range('E11').Value = '  +0.13%  '

$ws.Range('D12').Value = '7.55'
$ws.Range('E12').Value = '  -1.04%  '

$ws.Range('E13').Value = '  -1.90%  '

$ws.Range('D14').Value = '2.911.66'

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '15.17'
$ws.Range('E15').Value = '  +1.13%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.454.99'
$ws.Range('E16').Value = '  -2.61%  '

$ws.Range('D17').Value = '0.845'
$ws.Range('E17').Value = '  -1.75%  '

$ws.Range('D18').Value = '42.882.29'
$ws.Range('E18').Value = '  +1.10%  '

$ws.Range('D19').Value = '12.91'
$ws.Range('E19').Value = '  +0.52%  '

$ws.Range('D20').Value = '6.65'
$ws.Range('E20').Value = '  +3.46%  '

$ws.Range('D21').Value = '0.0₃0962'
$ws.Range('E21').Value = '  -0.77%  '

$ws.Range('D22').Value = '69.78'
$ws.Range('E22').Value = '  -1.63%  '

$ws.Range('D23').Value = '250.51'
$ws.Range('E23').Value = '  +0.03%  '

$ws.Range('D24').Value = '2.95'
$ws.Range('E24').Value = '  +2.05%  '

$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  +0.16%  '

$ws.Range('D26').Value = '26.71'
$ws.Range('E26').Value = '  -0.94%  '

$ws.Range('E27').Value = '  -0.06%  '

$ws.Range('E28').Value = '  +4.10%  '

$ws.Range('D29').Value = '39.99'
$ws.Range('E29').Value = '  +4.74%  '

$ws.Range('D30').Value = '10.22'
$ws.Range('E30').Value = '  -0.12%  '

$ws.Range('D31').Value = '5.98'
$ws.Range('E31').Value = '  +1.02%  '

$ws.Range('D32').Value = '155.42'
$ws.Range('E32').Value = '  +0.44%  '

$ws.Range('E33').Value = '  +2.71%  '

$ws.Range('D34').Value = '19.01'
$ws.Range('E34').Value = '  +2.24%  '

$ws.Range('D35').Value = '3.28'
$ws.Range('E35').Value = '  -0.93%  '

$ws.Range('D36').Value = '0.0786'
$ws.Range('E36').Value = '  +0.16%  '

$ws.Range('E37').Value = '  -0.56%  '

$ws.Range('E38').Value = '  -2.68%  '

$ws.Range('E39').Value = '  -0.27%  '

$ws.Range('D40').Value = '23.61'
$ws.Range('E40').Value = '  -2.40%  '

$ws.Range('E41').Value = '  +14.12%  '

$ws.Range('D42').Value = '0.0303'
$ws.Range('E42').Value = '  +1.64%  '

$ws.Range('E43').Value = '  +0.31%  '

$ws.Range('D44').Value = '3.76'
$ws.Range('E44').Value = '  -1.57%  '

$ws.Range('D45').Value = '3.28'
$ws.Range('E45').Value = '  -2.77%  '

$ws.Range('D46').Value = '2.017.42'
$ws.Range('E46').Value = '  -0.58%  '

$ws.Range('D47').Value = '85.43'
$ws.Range('E47').Value = '  +1.01%  '

$ws.Range('D48').Value = '8.77'
$ws.Range('E48').Value = '  -1.86%  '

$ws.Range('D49').Value = '2.767.16'
$ws.Range('E49').Value = '  +0.08%  '

$ws.Range('D50').Value = '73.28'
$ws.Range('E50').Value = '  +2.23%  '

$ws.Range('D51').Value = '102.21'
$ws.Range('E51').Value = '  +0.81%  '
